# Generate Report for Handoff
# Update the "Latest Handoff Date"/"Latest Handoff Datetime" for the
# f8d52f4d-2026-4b61-b734-e4ffb3e10d19 file row (row 5) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D is "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-17-15 03:17:47"

# zh-cn sheet: column E is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-15 03:17:39"

# de-de sheet: column E is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-15 03:17:47"
